$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.855.91'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.868.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.7009'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07783'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.05%  '
$ws.Range('E9').Value = '  +2.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07850'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.192'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.93%  '
$ws.Range('B13').Value = 'Litecoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '92.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.77%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.860.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6976'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.647'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.844.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.89%  '
$ws.Range('E18').Value = '  +1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.56'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.114.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.73%  '
$ws.Range('E21').Value = '  +1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.670'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1516'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.996'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.47'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.549'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.298'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.243'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.204'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05113'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7901'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.939'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.169'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.712'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.337.41'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01888'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.752'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9659'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.073'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +11.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '107.30'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9997'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000127'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.96%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.832'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.012.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.43%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.803'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.30%  '
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.046'
$ws.Range('D51').Style = 'Normal'
